$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Activity Diagram"
$ws.Range("M3").Value = "Activity diagram"
$ws.Range("I3").Value = "Use case diadram"

$ws.Range("I3").Select()
